# Scheduled-runner market data refresh for the Leve profit tables.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets to the
# latest Universalis market averages, matching upstream commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1721.4445
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1721.4445
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5164.333500000001
$ws.Range("N17").Value = -5500.333500000001
# Row 87
$ws.Range("H87").Value = 16285.714
$ws.Range("I87").Value = 8000
$ws.Range("J87").Value = 17666.666
$ws.Range("K87").Value = 8000
$ws.Range("L87").Value = 17666.666
$ws.Range("M87").Value = -6752
$ws.Range("N87").Value = -20162.666
# Row 90
$ws.Range("H90").Value = 16285.714
$ws.Range("I90").Value = 8000
$ws.Range("J90").Value = 17666.666
$ws.Range("K90").Value = 24000
$ws.Range("L90").Value = 52999.99800000001
$ws.Range("M90").Value = -17760
$ws.Range("N90").Value = -65479.99800000001
# Row 92
$ws.Range("H92").Value = 10083.625
$ws.Range("I92").Value = 140.66667
$ws.Range("J92").Value = 16049.4
$ws.Range("K92").Value = 140.66667
$ws.Range("L92").Value = 16049.4
$ws.Range("M92").Value = 1107.33333
$ws.Range("N92").Value = -18545.4
# Row 114
$ws.Range("H114").Value = 76000
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 76000
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 76000
$ws.Range("M114").ClearContents()
$ws.Range("N114").Value = -84678
# Row 137
$ws.Range("H137").Value = 7059.9355
$ws.Range("I137").Value = 2641.1538
$ws.Range("J137").Value = 10251.277
$ws.Range("K137").Value = 7923.4614
$ws.Range("L137").Value = 30753.831
$ws.Range("M137").Value = -5373.4614
# Row 138
$ws.Range("H138").Value = 9452.701999999999
$ws.Range("I138").Value = 5462.636
$ws.Range("J138").Value = 9981.505999999999
$ws.Range("K138").Value = 16387.908
$ws.Range("L138").Value = 29944.518
$ws.Range("M138").Value = -11247.908
$ws.Range("N138").Value = -40224.518
# Row 141
$ws.Range("H141").Value = 3428.9167
$ws.Range("I141").Value = 3536.5
$ws.Range("J141").Value = 3352.0715
$ws.Range("K141").Value = 10609.5
$ws.Range("L141").Value = 10056.2145
$ws.Range("M141").Value = -5429.5

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value = 11333.333
$ws.Range("I22").Value = 15500
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 15500
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -15201
$ws.Range("N22").Value = -3598
# Row 32
$ws.Range("H32").Value = 23390.816
$ws.Range("I32").Value = 14449.07
$ws.Range("J32").Value = 59796.5
$ws.Range("K32").Value = 14449.07
$ws.Range("L32").Value = 59796.5
$ws.Range("M32").Value = -14162.07
$ws.Range("N32").Value = -60370.5
# Row 61
$ws.Range("H61").Value = 6298.6553
$ws.Range("I61").Value = 5718.067
$ws.Range("J61").Value = 6920.7144
$ws.Range("K61").Value = 5718.067
$ws.Range("L61").Value = 6920.7144
$ws.Range("M61").Value = -5506.067
$ws.Range("N61").Value = -7344.7144
# Row 74
$ws.Range("H74").Value = 1840.3077
$ws.Range("I74").Value = 1674.5
$ws.Range("J74").Value = 1914
$ws.Range("K74").Value = 1674.5
$ws.Range("L74").Value = 1914
$ws.Range("M74").Value = -800.5
# Row 77
$ws.Range("H77").Value = 1840.3077
$ws.Range("I77").Value = 1674.5
$ws.Range("J77").Value = 1914
$ws.Range("K77").Value = 8372.5
$ws.Range("L77").Value = 9570
$ws.Range("M77").Value = -4004.5
# Row 132
$ws.Range("H132").Value = 4323.2593
$ws.Range("I132").Value = 4996.467
$ws.Range("J132").Value = 3481.75
$ws.Range("K132").Value = 14989.401
$ws.Range("L132").Value = 10445.25
$ws.Range("M132").Value = -12459.401
# Row 136
$ws.Range("H136").Value = 6298.6553
$ws.Range("I136").Value = 5718.067
$ws.Range("J136").Value = 6920.7144
$ws.Range("K136").Value = 17154.201
$ws.Range("L136").Value = 20762.1432
$ws.Range("M136").Value = -14604.201
$ws.Range("N136").Value = -25862.1432

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2574
$ws.Range("I94").Value = 669.5
$ws.Range("J94").Value = 3526.25
$ws.Range("K94").Value = 669.5
$ws.Range("L94").Value = 3526.25
$ws.Range("M94").Value = -218.5
$ws.Range("N94").Value = -4428.25
# Row 97
$ws.Range("H97").Value = 2138.75
$ws.Range("I97").Value = 2138.75
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2138.75
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1147.75
# Row 134
$ws.Range("H134").Value = 4761.2905
$ws.Range("I134").Value = 3178.0667
$ws.Range("J134").Value = 6245.5625
$ws.Range("K134").Value = 9534.2001
$ws.Range("L134").Value = 18736.6875
$ws.Range("M134").Value = -6999.2001

$ws = $wb.Worksheets.Item("CRP")
# Row 64
$ws.Range("H64").Value = 50000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 50000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
# Row 67
$ws.Range("H67").Value = 50000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 50000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
# Row 94
$ws.Range("H94").Value = 2266.7778
$ws.Range("I94").Value = 1898
$ws.Range("J94").Value = 3004.3333
$ws.Range("K94").Value = 1898
$ws.Range("L94").Value = 3004.3333
$ws.Range("M94").Value = -1447
$ws.Range("N94").Value = -3906.3333
# Row 122
$ws.Range("H122").Value = 3142.258
$ws.Range("I122").Value = 2829
$ws.Range("J122").Value = 3576
$ws.Range("K122").Value = 8487
$ws.Range("L122").Value = 10728
$ws.Range("M122").Value = -6037
$ws.Range("N122").Value = -15628
# Row 132
$ws.Range("H132").Value = 4785.4
$ws.Range("I132").Value = 4238
$ws.Range("J132").Value = 6975
$ws.Range("K132").Value = 12714
$ws.Range("L132").Value = 20925
$ws.Range("M132").Value = -10184
$ws.Range("N132").Value = -25985
# Row 134
$ws.Range("H134").Value = 2077.6956
$ws.Range("I134").Value = 1881.7059
$ws.Range("J134").Value = 2633
$ws.Range("K134").Value = 5645.1177
$ws.Range("L134").Value = 7899
$ws.Range("M134").Value = -3110.1177

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 4251794
$ws.Range("I4").Value = 5001828
$ws.Range("J4").Value = 1599
$ws.Range("K4").Value = 15005484
$ws.Range("L4").Value = 4797
$ws.Range("M4").Value = -15005372
$ws.Range("N4").Value = -5021
# Row 23
$ws.Range("H23").Value = 156.3
$ws.Range("I23").Value = 83
$ws.Range("J23").Value = 266.25
$ws.Range("K23").Value = 249
$ws.Range("L23").Value = 798.75
$ws.Range("M23").Value = -14
$ws.Range("N23").Value = -1268.75
# Row 33
$ws.Range("H33").Value = 2373.889
$ws.Range("I33").Value = 121.666664
$ws.Range("J33").Value = 3500
$ws.Range("K33").Value = 729.999984
$ws.Range("L33").Value = 21000
$ws.Range("M33").Value = -446.999984
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
# Row 55
$ws.Range("H55").Value = 50325
$ws.Range("I55").Value = 650
$ws.Range("J55").Value = 100000
$ws.Range("K55").Value = 1950
$ws.Range("L55").Value = 300000
$ws.Range("M55").Value = -1773
# Row 103
$ws.Range("H103").Value = 524.75
$ws.Range("I103").Value = 393
$ws.Range("J103").Value = 920
$ws.Range("K103").Value = 1179
$ws.Range("L103").Value = 2760
$ws.Range("M103").Value = -300
$ws.Range("N103").Value = -4518
# Row 113
$ws.Range("H113").Value = 1200.6666
$ws.Range("I113").Value = 100
$ws.Range("J113").Value = 1751
$ws.Range("K113").Value = 300
$ws.Range("L113").Value = 5253
$ws.Range("M113").Value = 1870
$ws.Range("N113").Value = -9593
# Row 118
$ws.Range("H118").Value = 838.1429000000001
$ws.Range("I118").Value = 838.1429000000001
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 2514.4287
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -1271.4287
# Row 121
$ws.Range("H121").Value = 2192.4443
$ws.Range("I121").Value = 3075.5
$ws.Range("J121").Value = 1940.1428
$ws.Range("K121").Value = 9226.5
$ws.Range("L121").Value = 5820.428400000001
$ws.Range("M121").Value = -7916.5
# Row 122
$ws.Range("H122").Value = 2518.8718
$ws.Range("I122").Value = 2299.5
$ws.Range("J122").Value = 2530.7297
$ws.Range("K122").Value = 20695.5
$ws.Range("L122").Value = 22776.5673
$ws.Range("M122").Value = -18245.5
$ws.Range("N122").Value = -27676.5673
# Row 131
$ws.Range("H131").Value = 26656.38
$ws.Range("I131").Value = 186202
$ws.Range("J131").Value = 7120.1836
$ws.Range("K131").Value = 558606
$ws.Range("L131").Value = 21360.5508
$ws.Range("M131").Value = -553566
$ws.Range("N131").Value = -31440.5508
# Row 132
$ws.Range("H132").Value = 2592.7646
$ws.Range("I132").Value = 1501.125
$ws.Range("J132").Value = 2928.6538
$ws.Range("K132").Value = 13510.125
$ws.Range("L132").Value = 26357.8842
$ws.Range("M132").Value = -10980.125
$ws.Range("N132").Value = -31417.8842
# Row 134
$ws.Range("H134").Value = 5107.4546
$ws.Range("I134").Value = 4152.3335
$ws.Range("J134").Value = 5465.625
$ws.Range("K134").Value = 12457.0005
$ws.Range("L134").Value = 16396.875
$ws.Range("M134").Value = -7387.000499999998

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 12065.833
$ws.Range("I43").Value = 4960
$ws.Range("J43").Value = 17141.428
$ws.Range("K43").Value = 4960
$ws.Range("L43").Value = 17141.428
$ws.Range("M43").Value = -4809
$ws.Range("N43").Value = -17443.428
# Row 93
$ws.Range("H93").Value = 37736.625
# Row 113
$ws.Range("H113").Value = 1600.7778
$ws.Range("I113").Value = 1644.125
$ws.Range("J113").Value = 1254
$ws.Range("K113").Value = 1644.125
$ws.Range("L113").Value = 1254
$ws.Range("M113").Value = 525.875
$ws.Range("N113").Value = -5594

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5041
$ws.Range("I40").Value = 5222.2144
$ws.Range("J40").Value = 4759.1113
$ws.Range("K40").Value = 5222.2144
$ws.Range("L40").Value = 4759.1113
$ws.Range("M40").Value = -5086.2144
$ws.Range("N40").Value = -5031.1113
# Row 55
$ws.Range("H55").Value = 924.6667
$ws.Range("I55").Value = 924.6667
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 924.6667
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -751.6667
# Row 132
$ws.Range("H132").Value = 4136.263
$ws.Range("I132").Value = 4313.7144
$ws.Range("J132").Value = 4032.75
$ws.Range("K132").Value = 12941.1432
$ws.Range("L132").Value = 12098.25
$ws.Range("M132").Value = -10411.1432
$ws.Range("N132").Value = -17158.25
# Row 136
$ws.Range("H136").Value = 3960.0833
$ws.Range("I136").Value = 4228.467
$ws.Range("J136").Value = 3512.7778
$ws.Range("K136").Value = 12685.401
$ws.Range("L136").Value = 10538.3334
$ws.Range("M136").Value = -10135.401

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
# Row 62
$ws.Range("H62").Value = 31769.9
$ws.Range("I62").Value = 12966.333
$ws.Range("J62").Value = 39828.57
$ws.Range("K62").Value = 12966.333
$ws.Range("L62").Value = 39828.57
$ws.Range("M62").Value = -12342.333
$ws.Range("N62").Value = -41076.57
# Row 65
$ws.Range("H65").Value = 31769.9
$ws.Range("I65").Value = 12966.333
$ws.Range("J65").Value = 39828.57
$ws.Range("K65").Value = 64831.665
$ws.Range("L65").Value = 199142.85
$ws.Range("M65").Value = -61711.665
$ws.Range("N65").Value = -205382.85
# Row 81
$ws.Range("H81").Value = 2141
$ws.Range("I81").Value = 2141
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 4282
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -3221
# Row 84
$ws.Range("H84").Value = 2141
$ws.Range("I84").Value = 2141
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 21410
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -16106
# Row 104
$ws.Range("H104").Value = 22058.375
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 22058.375
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 22058.375
$ws.Range("N104").Value = -29046.375
# Row 113
$ws.Range("H113").Value = 855.3043
$ws.Range("I113").Value = 758.8
$ws.Range("J113").Value = 1036.25
$ws.Range("K113").Value = 2276.4
$ws.Range("L113").Value = 3108.75
$ws.Range("M113").Value = -106.3999999999996

